$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Octubre de 2020 a las 16:28"

# Re-sort: Birmania now comes before Australia in the data table.
# Row 82 used to be Australia, row 83 used to be Birmania; the country
# names swap places and the statistics follow the (now updated) country.
$ws.Range("A82").Value = "Birmania"
$ws.Range("A83").Value = "Australia"

# Updated statistics per row (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes)
$data = @{
    4   = @(7949913, 4408, 5090255, 2640354, 0, 22, 219304)
    17  = @(481371, 1776, 453352, 14701, 0, 46, 13318)
    18  = @(402330, 2206, 336157, 56321, 0, 62, 9852)
    25  = @(324326, 873, 273500, 41129, 0, 6, 9697)
    58  = @(62151, 389, 44357, 16333, 0, 3, 1461)
    59  = @(61017, 241, 57967, 2545, 0, 3, 505)
    72  = @(41982, 230, 39357, 2016, 0, 1, 609)
    74  = @(41546, 388, 32000, 8780, 0, 6, 766)
    77  = @(34787, 102, 31536, 2488, 0, 1, 763)
    82  = @(27974, 1910, 9742, 17586, 0, 48, 646)
    83  = @(27265, 21, 24998, 1369, 0, 1, 898)
    95  = @(15506, 40, 11863, 3368, 0, 0, 275)
    102 = @(11936, 45, 9817, 1991, 0, 0, 128)
    145 = @(3526, 66, 2499, 1017, 0, 0, 10)
    177 = @(524, 7, 472, 51, 0, 0, 1)
    196 = @(144, 0, 121, 22, 0, 0, 1)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
    $ws.Range("G$row").Value = $vals[5]
    $ws.Range("H$row").Value = $vals[6]
}
